# Fixing the big mistake
# Update mean (row 3) and community (column D) statistics that were
# recalculated, plus the dependent Sums / % energy sector figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (mean): Total (B) and Community (D) columns
$ws.Range("B3").Value = 29064.13089906011
$ws.Range("D3").Value = 1258.796563160196

# Row 4 (std): Total (B) and Community (D) columns
$ws.Range("B4").Value = 12245.21117573238
$ws.Range("D4").Value = 693.7101388976731

# Row 5 (min): Total (B) and Community (D) columns
$ws.Range("B5").Value = 9662.851936986302
$ws.Range("D5").Value = 166.8318904109584

# Row 6 (25%): Total (B) and Community (D) columns
$ws.Range("B6").Value = 19950.67174452056
$ws.Range("D6").Value = 570.2154821917791

# Row 7 (50%): Total (B) and Community (D) columns
$ws.Range("B7").Value = 24702.6474616439
$ws.Range("D7").Value = 1111.802752054795

# Row 8 (75%): Total (B) and Community (D) columns
$ws.Range("B8").Value = 40358.65445273992
$ws.Range("D8").Value = 1996.000010958904

# Row 9 (max): Total (B) and Community (D) columns
$ws.Range("B9").Value = 53836.45736438355
$ws.Range("D9").Value = 2054.767134246576

# Row 10 (Total sums): F column
$ws.Range("F10").Value = 41852348.49464662

# Row 11 (Residential sums): % energy sector (G)
$ws.Range("G11").Value = 0.8208412111448016

# Row 12 (Community sums): F and % energy sector (G)
$ws.Range("F12").Value = 1812667.050950684
$ws.Range("G12").Value = 0.0433109996487424

# Row 13 (IGA sums): % energy sector (G)
$ws.Range("G13").Value = 0.135847789206456
